$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.234.79"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.590.07"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.05"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.246"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.24"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.812.86"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.588.76"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.13"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "26.233.20"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.08"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.98"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.98"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.01"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.16"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "1.408.92"
$ws.Range("E33").Value = "  +8.63%  "
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.586"
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.86"
$ws.Range("E40").Value = "  +4.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.952"
$ws.Range("E42").Value = "  -14.27%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.769"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "1.724.62"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.11"
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.33"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0501"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.22%  "
